{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the \"Minor edit to user stories\" changes described by the diff.\n//\n// Strategy: use body.search() to locate the exact old substrings (which\n// may straddle several <w:r> runs / proofErr markers) and replace them via\n// Range.insertText(..., Word.InsertLocation.replace) \u2014 this collapses the\n// matched range to a single new run with the replacement text, mirroring\n// the textual result of the diff. Then remove the one extraneous empty\n// paragraph that the diff deletes.\n\nasync function replaceOnce(body, searchText, replacement, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) Register acceptance criteria: reword register-button sentence and\n//    add the \"After filling in...\" sentence; \"password and\" -> \"password, and\".\nawait replaceOnce(\n  body,\n  \"clicks a register button and text boxes appear to enter email, username, password and city. Then the user can get access to the dashboard.\",\n  \"clicks on the register button and text boxes appear to enter email, username, password, and city. After filling in the information and creating an account, the user gets access to the dashboard.\"\n);\n\n// 2) Log in acceptance criteria: \"filling their\" -> \"filling in their\".\nawait replaceOnce(\n  body,\n  \"Users log in by filling their username and password\",\n  \"Users log in by filling in their username and password\"\n);\n\n// 3) Forgot username/password acceptance criteria: add a sentence about the\n//    new credentials being emailed.\nawait replaceOnce(\n  body,\n  \"Users click \\u201Cforgot username\\u201D or \\u201Cforgot password\\u201D button and fill in their emails. Then they can access to their accounts again.\",\n  \"Users click \\u201Cforgot username\\u201D or \\u201Cforgot password\\u201D button and fill in their emails. New username and/or password would be sent to their email and they can access to their accounts again.\"\n);\n\n// 4) Post-comments user story: the diff merges a run split (\"I ca\" / \"n\n//    share...\") into one run. No visible text change, but normalize it\n//    anyway via a no-op replace so the run boundary collapses like the\n//    canonical XML.\nawait replaceOnce(\n  body,\n  \"As a user, I want to be able to post comments on graphs, so that I can share my thought to other users.\",\n  \"As a user, I want to be able to post comments on graphs, so that I can share my thought to other users.\"\n);\n\n// 4b) Remove the extraneous empty (ind left=720) paragraph that sits\n//     between the empty paragraph and the \"Acceptance criteria\" paragraph\n//     for the comments user story.\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"\" && i > 0 && i + 1 < paras.items.length) {\n    const prevText = paras.items[i - 1].text;\n    const nextText = paras.items[i + 1].text;\n    if (\n      prevText === \"\" &&\n      nextText.indexOf(\"Users click \\u201Ccomment\\u201D button\") !== -1\n    ) {\n      paras.items[i].delete();\n      await context.sync();\n      break;\n    }\n  }\n}\n\n// 5) Comment acceptance criteria: add \", and submit\" before the period.\nawait replaceOnce(\n  body,\n  \"Users click \\u201Ccomment\\u201D button and type in their comment in the text box.\",\n  \"Users click \\u201Ccomment\\u201D button, type in their comment in the text box, and submit.\"\n);\n\n// 6) Delete-comments user story: \"delete comments\" -> \"delete other users' comments\".\nawait replaceOnce(\n  body,\n  \"As an admin, I want to be able to delete comments, so that I can get rid of any inappropriate comments.\",\n  \"As an admin, I want to be able to delete other users\\u2019 comments, so that I can get rid of any inappropriate comments.\"\n);\n\n// 7) Delete acceptance criteria: \"deleted from the comment.\" -> \"deleted from the graphs.\"\nawait replaceOnce(\n  body,\n  \"Admin clicks the button \\u201Cdelete\\u201D on a comment, and the comment is deleted from the comment.\",\n  \"Admin clicks the button \\u201Cdelete\\u201D on a comment, and the comment is deleted from the graphs.\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the \"Minor edit to user stories\" changes described by the diff.\n#\n# Strategy: use Range.Find.Execute(..., Replace:=wdReplaceOne) against the\n# whole document content to locate the exact old sentences (which may\n# straddle several runs / proofErr markers) and replace them in one shot.\n# Then remove the one extraneous empty paragraph that the diff deletes.\n\n$d = $word.ActiveDocument\n\nfunction ReplaceText($findText, $replaceText) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $found = $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"ReplaceText: could not find: $findText\"\n    }\n}\n\n# 1) Register acceptance criteria: reword register-button sentence and\n#    add the \"After filling in...\" sentence; \"password and\" -> \"password, and\".\nReplaceText \"clicks a register button and text boxes appear to enter email, username, password and city. Then the user can get access to the dashboard.\" \"clicks on the register button and text boxes appear to enter email, username, password, and city. After filling in the information and creating an account, the user gets access to the dashboard.\"\n\n# 2) Log in acceptance criteria: \"filling their\" -> \"filling in their\".\nReplaceText \"Users log in by filling their username and password\" \"Users log in by filling in their username and password\"\n\n# 3) Forgot username/password acceptance criteria: add a sentence about the\n#    new credentials being emailed.\nReplaceText \"Users click \u201cforgot username\u201d or \u201cforgot password\u201d button and fill in their emails. Then they can access to their accounts again.\" \"Users click \u201cforgot username\u201d or \u201cforgot password\u201d button and fill in their emails. New username and/or password would be sent to their email and they can access to their accounts again.\"\n\n# 4) Post-comments user story: the diff merges a run split (\"I ca\" / \"n\n#    share...\") into one run. No visible text change, but normalize it\n#    anyway via a no-op replace so the run boundary collapses like the\n#    canonical XML.\nReplaceText \"As a user, I want to be able to post comments on graphs, so that I can share my thought to other users.\" \"As a user, I want to be able to post comments on graphs, so that I can share my thought to other users.\"\n\n# 4b) Remove the extraneous empty (ind left=720) paragraph that sits\n#     between the empty paragraph and the \"Acceptance criteria\" paragraph\n#     for the comments user story.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -eq \"`r\" -and $i -gt 1 -and $i -lt $count) {\n        $prevText = $d.Paragraphs($i - 1).Range.Text\n        $nextText = $d.Paragraphs($i + 1).Range.Text\n        if ($prevText -eq \"`r\" -and $nextText -like \"*Users click `\u201ccomment`\u201d button*\") {\n            $p.Range.Delete()\n            break\n        }\n    }\n}\n\n# 5) Comment acceptance criteria: add \", and submit\" before the period.\nReplaceText \"Users click \u201ccomment\u201d button and type in their comment in the text box.\" \"Users click \u201ccomment\u201d button, type in their comment in the text box, and submit.\"\n\n# 6) Delete-comments user story: \"delete comments\" -> \"delete other users' comments\".\nReplaceText \"As an admin, I want to be able to delete comments, so that I can get rid of any inappropriate comments.\" \"As an admin, I want to be able to delete other users\u2019 comments, so that I can get rid of any inappropriate comments.\"\n\n# 7) Delete acceptance criteria: \"deleted from the comment.\" -> \"deleted from the graphs.\"\nReplaceText \"Admin clicks the button \u201cdelete\u201d on a comment, and the comment is deleted from the comment.\" \"Admin clicks the button \u201cdelete\u201d on a comment, and the comment is deleted from the graphs.\"\n"}
